$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 102 <= old row 103
$ws.Range("F102").Value = 'Albirex Niigata'
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 'Kashiwa Reysol'
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 2.22
$ws.Range("K102").Value = '03/05/2023 12:12'
$ws.Range("L102").Value = 2.73
$ws.Range("M102").Value = '07/05/2023 06:52'
$ws.Range("N102").Value = 3.39
$ws.Range("O102").Value = '03/05/2023 12:12'
$ws.Range("P102").Value = 3.23
$ws.Range("Q102").Value = '07/05/2023 06:41'
$ws.Range("R102").Value = 3.48
$ws.Range("S102").Value = '03/05/2023 12:12'
$ws.Range("T102").Value = 2.86
$ws.Range("U102").Value = '07/05/2023 06:52'
$ws.Range("V102").Value = 'https://www.betexplorer.com/football/japan/j1-league/albirex-niigata-kashiwa-reysol/6XPoZcmN/'

# row 103 <= old row 105
$ws.Range("F103").Value = 'Vissel Kobe'
$ws.Range("G103").Value = 3
$ws.Range("H103").Value = 'Yokohama FC'
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 1.39
$ws.Range("K103").Value = '03/05/2023 12:12'
$ws.Range("L103").Value = 1.38
$ws.Range("M103").Value = '07/05/2023 06:41'
$ws.Range("N103").Value = 5.05
$ws.Range("O103").Value = '03/05/2023 12:12'
$ws.Range("P103").Value = 5.27
$ws.Range("Q103").Value = '07/05/2023 06:41'
$ws.Range("R103").Value = 7.97
$ws.Range("S103").Value = '03/05/2023 12:12'
$ws.Range("T103").Value = 8.31
$ws.Range("U103").Value = '07/05/2023 06:41'
$ws.Range("V103").Value = 'https://www.betexplorer.com/football/japan/j1-league/vissel-kobe-yokohama-fc/z3rATyBp/'

# row 105 <= old row 102
$ws.Range("F105").Value = 'Yokohama F. Marinos'
$ws.Range("G105").Value = 4
$ws.Range("H105").Value = 'Kyoto'
$ws.Range("I105").Value = 1
$ws.Range("J105").Value = 1.48
$ws.Range("K105").Value = '03/05/2023 08:12'
$ws.Range("L105").Value = 1.55
$ws.Range("M105").Value = '07/05/2023 06:51'
$ws.Range("N105").Value = 5.02
$ws.Range("O105").Value = '03/05/2023 08:12'
$ws.Range("P105").Value = 4.95
$ws.Range("Q105").Value = '07/05/2023 06:52'
$ws.Range("R105").Value = 5.84
$ws.Range("S105").Value = '03/05/2023 08:12'
$ws.Range("T105").Value = 5.46
$ws.Range("U105").Value = '07/05/2023 06:52'
$ws.Range("V105").Value = 'https://www.betexplorer.com/football/japan/j1-league/yokohama-f-marinos-kyoto/CAQszvYG/'

# row 240 <= old row 242
$ws.Range("F240").Value = 'Sanfrecce Hiroshima'
$ws.Range("G240").Value = 2
$ws.Range("H240").Value = 'Vissel Kobe'
$ws.Range("I240").Value = 0
$ws.Range("J240").Value = 2.04
$ws.Range("K240").Value = '04/09/2023 08:42'
$ws.Range("L240").Value = 2.06
$ws.Range("M240").Value = '16/09/2023 11:56'
$ws.Range("N240").Value = 3.58
$ws.Range("O240").Value = '04/09/2023 08:42'
$ws.Range("P240").Value = 3.63
$ws.Range("Q240").Value = '16/09/2023 11:58'
$ws.Range("R240").Value = 3.8
$ws.Range("S240").Value = '04/09/2023 08:42'
$ws.Range("T240").Value = 3.73
$ws.Range("U240").Value = '16/09/2023 11:56'
$ws.Range("V240").Value = 'https://www.betexplorer.com/football/japan/j1-league/sanfrecce-hiroshima-vissel-kobe/jex48het/'

# row 241 <= old row 240
$ws.Range("F241").Value = 'Avispa Fukuoka'
$ws.Range("G241").Value = 1
$ws.Range("H241").Value = 'Nagoya Grampus'
$ws.Range("I241").Value = 0
$ws.Range("J241").Value = 2.94
$ws.Range("K241").Value = '04/09/2023 08:42'
$ws.Range("L241").Value = 3.73
$ws.Range("M241").Value = '16/09/2023 11:59'
$ws.Range("N241").Value = 3.2
$ws.Range("O241").Value = '04/09/2023 08:42'
$ws.Range("P241").Value = 2.98
$ws.Range("Q241").Value = '16/09/2023 11:57'
$ws.Range("R241").Value = 2.64
$ws.Range("S241").Value = '04/09/2023 08:42'
$ws.Range("T241").Value = 2.35
$ws.Range("U241").Value = '16/09/2023 11:59'
$ws.Range("V241").Value = 'https://www.betexplorer.com/football/japan/j1-league/avispa-fukuoka-nagoya-grampus/zqOP3U9B/'

# row 242 <= old row 241
$ws.Range("F242").Value = 'Hokkaido Consadole Sapporo'
$ws.Range("G242").Value = 0
$ws.Range("H242").Value = 'Shonan Bellmare'
$ws.Range("I242").Value = 1
$ws.Range("J242").Value = 1.67
$ws.Range("K242").Value = '04/09/2023 08:42'
$ws.Range("L242").Value = 1.65
$ws.Range("M242").Value = '16/09/2023 11:48'
$ws.Range("N242").Value = 4.34
$ws.Range("O242").Value = '04/09/2023 08:42'
$ws.Range("P242").Value = 4.49
$ws.Range("Q242").Value = '16/09/2023 11:58'
$ws.Range("R242").Value = 4.87
$ws.Range("S242").Value = '04/09/2023 08:42'
$ws.Range("T242").Value = 4.92
$ws.Range("U242").Value = '16/09/2023 11:55'
$ws.Range("V242").Value = 'https://www.betexplorer.com/football/japan/j1-league/hokkaido-consadole-sapporo-shonan-bellmare/YJR87CAn/'
# new row 251
$ws.Cells.Item(250,1).Copy()
$ws.Cells.Item(251,1).PasteSpecial(-4122)
$ws.Cells.Item(250,5).Copy()
$ws.Cells.Item(251,5).PasteSpecial(-4122)
$ws.Range("A251").Value = 250
$ws.Range("B251").Value = 'japan'
$ws.Range("C251").Value = 'j1-league'
$ws.Range("D251").Value = '2023'
$ws.Range("E251").Value = 45193.33333333334
$ws.Range("F251").Value = 'Kashima Antlers'
$ws.Range("G251").Value = 1
$ws.Range("H251").Value = 'Yokohama F. Marinos'
$ws.Range("I251").Value = 2
$ws.Range("J251").Value = 2.24
$ws.Range("K251").Value = '17/09/2023 07:12'
$ws.Range("L251").Value = 2.35
$ws.Range("M251").Value = '24/09/2023 07:55'
$ws.Range("N251").Value = 3.72
$ws.Range("O251").Value = '17/09/2023 07:12'
$ws.Range("P251").Value = 3.74
$ws.Range("Q251").Value = '24/09/2023 07:55'
$ws.Range("R251").Value = 3.09
$ws.Range("S251").Value = '17/09/2023 07:12'
$ws.Range("T251").Value = 2.99
$ws.Range("U251").Value = '24/09/2023 07:55'
$ws.Range("V251").Value = 'https://www.betexplorer.com/football/japan/j1-league/kashima-antlers-yokohama-f-marinos/8dMX18vO/'

# new row 252
$ws.Cells.Item(250,1).Copy()
$ws.Cells.Item(252,1).PasteSpecial(-4122)
$ws.Cells.Item(250,5).Copy()
$ws.Cells.Item(252,5).PasteSpecial(-4122)
$ws.Range("A252").Value = 251
$ws.Range("B252").Value = 'japan'
$ws.Range("C252").Value = 'j1-league'
$ws.Range("D252").Value = '2023'
$ws.Range("E252").Value = 45193.375
$ws.Range("F252").Value = 'Shonan Bellmare'
$ws.Range("G252").Value = 0
$ws.Range("H252").Value = 'Kawasaki Frontale'
$ws.Range("I252").Value = 2
$ws.Range("J252").Value = 3.07
$ws.Range("K252").Value = '17/09/2023 08:12'
$ws.Range("L252").Value = 4.36
$ws.Range("M252").Value = '24/09/2023 08:59'
$ws.Range("N252").Value = 3.61
$ws.Range("O252").Value = '17/09/2023 08:12'
$ws.Range("P252").Value = 3.94
$ws.Range("Q252").Value = '24/09/2023 08:55'
$ws.Range("R252").Value = 2.29
$ws.Range("S252").Value = '17/09/2023 08:12'
$ws.Range("T252").Value = 1.83
$ws.Range("U252").Value = '24/09/2023 08:55'
$ws.Range("V252").Value = 'https://www.betexplorer.com/football/japan/j1-league/shonan-bellmare-kawasaki-frontale/nNirMBPu/'

# new row 253
$ws.Cells.Item(250,1).Copy()
$ws.Cells.Item(253,1).PasteSpecial(-4122)
$ws.Cells.Item(250,5).Copy()
$ws.Cells.Item(253,5).PasteSpecial(-4122)
$ws.Range("A253").Value = 252
$ws.Range("B253").Value = 'japan'
$ws.Range("C253").Value = 'j1-league'
$ws.Range("D253").Value = '2023'
$ws.Range("E253").Value = 45193.41666666666
$ws.Range("F253").Value = 'Gamba Osaka'
$ws.Range("G253").Value = 1
$ws.Range("H253").Value = 'Urawa Reds'
$ws.Range("I253").Value = 3
$ws.Range("J253").Value = 2.83
$ws.Range("K253").Value = '17/09/2023 11:12'
$ws.Range("L253").Value = 3.44
$ws.Range("M253").Value = '24/09/2023 09:56'
$ws.Range("N253").Value = 3.33
$ws.Range("O253").Value = '17/09/2023 11:12'
$ws.Range("P253").Value = 3.37
$ws.Range("Q253").Value = '24/09/2023 09:56'
$ws.Range("R253").Value = 2.64
$ws.Range("S253").Value = '17/09/2023 11:12'
$ws.Range("T253").Value = 2.27
$ws.Range("U253").Value = '24/09/2023 09:56'
$ws.Range("V253").Value = 'https://www.betexplorer.com/football/japan/j1-league/gamba-osaka-urawa-reds/t4gbITO4/'
